$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 5, 95.4377950037229, 0, 0, 0, 143.1566925055843, 0.39, 0.4, 0.15, 0.35),
    @(5, 6, 95.10188558031456, 0, 0, 0, 142.6528283704718, 0.34, 0.5, 0.11, 0.34),
    @(6, 7, 100.6937895111742, 0, 0, 0, 151.0406842667613, 0.39, 0.45, 0.14, 0.34),
    @(7, 8, 122.9626083453766, 0, 0, 0, 184.4439125180649, 0.4, 0.44, 0.19, 0.35),
    @(8, 9, 217.7285845022828, 0.00555556, 0.00555556, 0.01111112, 326.5928767534242, 0.38, 0.41, 0.17, 0.33),
    @(9, 10, 340.6121553362682, 0.02222224, 0.02222224, 0.04444448, 510.9182330044023, 0.38, 0.47, 0.18, 0.39),
    @(10, 11, 428.3635523572899, 0.08611118000000001, 0.08611118000000001, 0.17222236, 642.5453285359349, 0.4, 0.45, 0.12, 0.4),
    @(11, 12, 498.7069374946091, 0.10000008, 0.10000008, 0.20000016, 748.0604062419137, 0.34, 0.42, 0.14, 0.3),
    @(12, 13, 508.3495138842192, 0.11944454, 0.11944454, 0.23888908, 762.5242708263288, 0.35, 0.48, 0.12, 0.36),
    @(13, 14, 425.5577207029365, 0.10277786, 0.10277786, 0.20555572, 638.3365810544047, 0.35, 0.45, 0.14, 0.38),
    @(14, 15, 334.7633794933903, 0.0416667, 0.0416667, 0.0833334, 502.1450692400854, 0.3, 0.4, 0.11, 0.37),
    @(15, 16, 286.7480913238156, 0.0138889, 0.0138889, 0.0277778, 430.1221369857234, 0.33, 0.49, 0.12, 0.36),
    @(16, 17, 262.4045378156197, 0, 0, 0, 393.6068067234295, 0.39, 0.46, 0.14, 0.39),
    @(17, 18, 309.1354664251031, 0, 0, 0, 463.7031996376547, 0.37, 0.48, 0.19, 0.34),
    @(18, 19, 384.636049180636, 0, 0, 0, 576.954073770954, 0.35, 0.4, 0.18, 0.39),
    @(19, 20, 433.5207499755035, 0, 0, 0, 650.2811249632553, 0.38, 0.41, 0.17, 0.33)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# New rows 13:17 need the same style as column A of the existing rows (A2:A12)
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
